$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last filled data row (309) down to the new rows (310-313)
$ws.Range("A309:C309").Copy()
$ws.Range("A310:C313").PasteSpecial(-4122)

$newRows = @(
    @("國際聯合科技股份有限公司", 16450332, "出口報單買方需加上地址"),
    @("譁泰精機股份有限公司", 97241217, "出口報單買方需加上地址"),
    @("為升電裝工業股份有限公司", 23219346, "出口報單買方需加上地址"),
    @("鼎茂光電股份有限公司", 54156182, "出口報單買方需加上地址")
)

$r = 310
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
